$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (UID 2): source well A1 -> A2, dest well A2 -> A1, volume 500 -> 125, reagent "DNA ligase buffer" -> "DNA ligase"
$ws.Range("D3").Value = "A2"
$ws.Range("G3").Value = "A1"
$ws.Range("H3").Value = 125
$ws.Range("I3").Value = "DNA ligase"

# Row 4 (UID 3): source well A2 -> A3, dest well stays A1, volume 125 -> 250, reagent "DNA ligase" -> "BsmBI (NEB)"
$ws.Range("D4").Value = "A3"
$ws.Range("H4").Value = 250
$ws.Range("I4").Value = "BsmBI (NEB)"

# Remove old rows 5, 6 and 7 (UID 4, 5, 6) - no longer needed
$ws.Range("A5:A7").EntireRow.Delete()

$ws.Range("A1:I4").Columns.AutoFit()
